$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Added all friends view with filter" — fill in previously-blank GitHub
# activity scores and the Friends section's "Display Top Friends with
# Images" score (C44's Total Score SUM formula recalculates automatically).
$ws.Range("C8").Value = 9
$ws.Range("C9").Value = 16
$ws.Range("C19").Value = 10

# Reflect the author's scrolled viewport / new active selection.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("H25").Select()
